$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 3 (pushes the existing row 3 - PREPROD/0420172008382 /55277 - down to row 4)
$ws.Rows("3:3").Insert()

# New row 3 holds the SISE payment validation data: claim number (B3) and
# payment order number (C3), as mentioned in the commit message.
$ws.Range("A3").Value = "PREPROD"
$ws.Range("C3").Value = "'55293  "
$ws.Range("B3").Value = "'0420172008486 "
